$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 49875.906
$ws.Range("J40").Value = 2543.889
$ws.Range("L40").Value = 2543.889
$ws.Range("N40").Value = -2893.889
$ws.Range("H64").Value = 52999.95
$ws.Range("I64").Value = 145428.86
$ws.Range("J64").Value = 3230.5386
$ws.Range("K64").Value = 145428.86
$ws.Range("L64").Value = 3230.5386
$ws.Range("M64").Value = -145180.86
$ws.Range("N64").Value = -3726.5386
$ws.Range("H67").Value = 52999.95
$ws.Range("I67").Value = 145428.86
$ws.Range("J67").Value = 3230.5386
$ws.Range("K67").Value = 145428.86
$ws.Range("L67").Value = 3230.5386
$ws.Range("M67").Value = -144570.86
$ws.Range("N67").Value = -4946.5386
$ws.Range("H112").Value = 1053.5
$ws.Range("J112").Value = 1053.5
$ws.Range("L112").Value = 3160.5
$ws.Range("N112").Value = -5376.5
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H135").Value = 677.7692
$ws.Range("I135").Value = 677.7692
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6099.922799999999
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3564.922799999999
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 1167.4509
$ws.Range("I137").Value = 957.7954999999999
$ws.Range("J137").Value = 2485.2856
$ws.Range("K137").Value = 2873.3865
$ws.Range("L137").Value = 7455.8568
$ws.Range("M137").Value = -323.3864999999996
$ws.Range("N137").Value = -12555.8568
$ws.Range("H138").Value = 2733.243
$ws.Range("I138").Value = 1426.1538
$ws.Range("J138").Value = 3505.6135
$ws.Range("K138").Value = 4278.4614
$ws.Range("L138").Value = 10516.8405
$ws.Range("M138").Value = 861.5385999999999
$ws.Range("N138").Value = -20796.8405

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 78323.69500000001
$ws.Range("J45").Value = 2352
$ws.Range("L45").Value = 2352
$ws.Range("N45").Value = -3106
$ws.Range("H61").Value = 1883.0392
$ws.Range("I61").Value = 1136.9642
$ws.Range("J61").Value = 2791.3044
$ws.Range("K61").Value = 1136.9642
$ws.Range("L61").Value = 2791.3044
$ws.Range("M61").Value = -924.9641999999999
$ws.Range("N61").Value = -3215.3044
$ws.Range("H74").Value = 1321.8214
$ws.Range("I74").Value = 604.4211
$ws.Range("J74").Value = 2836.3333
$ws.Range("K74").Value = 604.4211
$ws.Range("L74").Value = 2836.3333
$ws.Range("M74").Value = 269.5789
$ws.Range("N74").Value = -4584.3333
$ws.Range("H77").Value = 1321.8214
$ws.Range("I77").Value = 604.4211
$ws.Range("J77").Value = 2836.3333
$ws.Range("K77").Value = 3022.1055
$ws.Range("L77").Value = 14181.6665
$ws.Range("M77").Value = 1345.8945
$ws.Range("N77").Value = -22917.6665
$ws.Range("H136").Value = 1883.0392
$ws.Range("I136").Value = 1136.9642
$ws.Range("J136").Value = 2791.3044
$ws.Range("K136").Value = 3410.8926
$ws.Range("L136").Value = 8373.913199999999
$ws.Range("M136").Value = -860.8925999999997
$ws.Range("N136").Value = -13473.9132

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1835.4286
$ws.Range("I134").Value = 1765.9215
$ws.Range("J134").Value = 2544.4
$ws.Range("K134").Value = 5297.764499999999
$ws.Range("L134").Value = 7633.200000000001
$ws.Range("M134").Value = -2762.764499999999
$ws.Range("N134").Value = -12703.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 13336
$ws.Range("I17").Value = 8
$ws.Range("J17").Value = 20000
$ws.Range("K17").Value = 8
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 166
$ws.Range("N17").Value = -20348
$ws.Range("H25").Value = 20000
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H41").Value = 8681
$ws.Range("I41").Value = 3762.5
$ws.Range("J41").Value = 11960
$ws.Range("K41").Value = 3762.5
$ws.Range("L41").Value = 11960
$ws.Range("M41").Value = -3334.5
$ws.Range("N41").Value = -12816
$ws.Range("H50").Value = 11205
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 11205
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 11205
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -12455
$ws.Range("H51").Value = 7841.8184
$ws.Range("J51").Value = 7841.8184
$ws.Range("L51").Value = 7841.8184
$ws.Range("N51").Value = -9313.8184
$ws.Range("H59").Value = 21607.5
$ws.Range("I59").Value = 8000
$ws.Range("J59").Value = 23551.428
$ws.Range("K59").Value = 8000
$ws.Range("L59").Value = 23551.428
$ws.Range("M59").Value = -6855
$ws.Range("N59").Value = -25841.428
$ws.Range("H60").Value = 9725.444
$ws.Range("I60").Value = 8820.799999999999
$ws.Range("J60").Value = 10073.385
$ws.Range("K60").Value = 8820.799999999999
$ws.Range("L60").Value = 10073.385
$ws.Range("M60").Value = -8309.799999999999
$ws.Range("N60").Value = -11095.385
$ws.Range("H61").Value = 7841.8184
$ws.Range("J61").Value = 7841.8184
$ws.Range("L61").Value = 7841.8184
$ws.Range("N61").Value = -8537.8184
$ws.Range("H68").Value = 16249
$ws.Range("J68").Value = 16249
$ws.Range("L68").Value = 16249
$ws.Range("N68").Value = -17747
$ws.Range("H71").Value = 16249
$ws.Range("J71").Value = 16249
$ws.Range("L71").Value = 48747
$ws.Range("N71").Value = -56235
$ws.Range("H108").Value = 27777
$ws.Range("J108").Value = 27777
$ws.Range("L108").Value = 27777
$ws.Range("N108").Value = -35457

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 553.125
$ws.Range("I34").Value = 112.5
$ws.Range("J34").Value = 993.75
$ws.Range("K34").Value = 337.5
$ws.Range("L34").Value = 2981.25
$ws.Range("M34").Value = -253.5
$ws.Range("N34").Value = -3149.25
$ws.Range("H122").Value = 671.55554
$ws.Range("I122").Value = 474.16666
$ws.Range("J122").Value = 1066.3334
$ws.Range("K122").Value = 4267.49994
$ws.Range("L122").Value = 9597.000599999999
$ws.Range("M122").Value = -1817.49994
$ws.Range("N122").Value = -14497.0006
$ws.Range("H131").Value = 836.23
$ws.Range("I131").Value = 530
$ws.Range("J131").Value = 862.8587
$ws.Range("K131").Value = 1590
$ws.Range("L131").Value = 2588.5761
$ws.Range("M131").Value = 3450
$ws.Range("N131").Value = -12668.5761
$ws.Range("H132").Value = 2117.2222
$ws.Range("J132").Value = 2436.4285
$ws.Range("L132").Value = 21927.8565
$ws.Range("N132").Value = -26987.8565

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 9719.777
$ws.Range("I46").Value = 5990
$ws.Range("J46").Value = 10785.429
$ws.Range("K46").Value = 5990
$ws.Range("L46").Value = 10785.429
$ws.Range("M46").Value = -5834
$ws.Range("N46").Value = -11097.429
$ws.Range("H49").Value = 43820
$ws.Range("J49").Value = 43820
$ws.Range("L49").Value = 43820
$ws.Range("N49").Value = -44188
$ws.Range("H57").Value = 11951.667
$ws.Range("I57").Value = 8027.5
$ws.Range("K57").Value = 8027.5
$ws.Range("M57").Value = -7207.5
$ws.Range("H80").Value = 125130376
$ws.Range("H83").Value = 125130376
$ws.Range("H113").Value = 1694.1765
$ws.Range("I113").Value = 1575.0834
$ws.Range("J113").Value = 1980
$ws.Range("K113").Value = 1575.0834
$ws.Range("L113").Value = 1980
$ws.Range("M113").Value = 594.9166
$ws.Range("N113").Value = -6320
$ws.Range("H122").Value = 3599.4
$ws.Range("I122").Value = 3599.4
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10798.2
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8348.200000000001
$ws.Range("N122").ClearContents()
$ws.Range("H124").Value = 48000
$ws.Range("J124").Value = 48000
$ws.Range("L124").Value = 48000
$ws.Range("N124").Value = -57820
$ws.Range("H126").Value = 3861.2856
$ws.Range("I126").Value = 3406
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 10218
$ws.Range("L126").Value = 14998.5
$ws.Range("M126").Value = -7748
$ws.Range("N126").Value = -19938.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6300739.5
$ws.Range("I16").Value = 9000529
$ws.Range("J16").Value = 1231.6666
$ws.Range("K16").Value = 9000529
$ws.Range("L16").Value = 1231.6666
$ws.Range("M16").Value = -9000359
$ws.Range("N16").Value = -1571.6666
$ws.Range("H22").Value = 895.7406999999999
$ws.Range("I22").Value = 768.4286
$ws.Range("J22").Value = 940.3
$ws.Range("K22").Value = 768.4286
$ws.Range("L22").Value = 940.3
$ws.Range("M22").Value = -473.4286
$ws.Range("N22").Value = -1530.3
$ws.Range("H27").Value = 895.7406999999999
$ws.Range("I27").Value = 768.4286
$ws.Range("J27").Value = 940.3
$ws.Range("K27").Value = 768.4286
$ws.Range("L27").Value = 940.3
$ws.Range("M27").Value = -661.4286
$ws.Range("N27").Value = -1154.3
$ws.Range("H42").Value = 9860
$ws.Range("J42").Value = 9860
$ws.Range("L42").Value = 9860
$ws.Range("N42").Value = -10986
$ws.Range("H46").Value = 779074.0600000001
$ws.Range("I46").Value = 413.33334
$ws.Range("J46").Value = 1012672.3
$ws.Range("K46").Value = 413.33334
$ws.Range("L46").Value = 1012672.3
$ws.Range("M46").Value = -225.33334
$ws.Range("N46").Value = -1013048.3
$ws.Range("H49").Value = 9860
$ws.Range("J49").Value = 9860
$ws.Range("L49").Value = 9860
$ws.Range("N49").Value = -10154
$ws.Range("H95").Value = 19000
$ws.Range("J95").Value = 19000
$ws.Range("L95").Value = 19000
$ws.Range("N95").Value = -24492
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 624.6
$ws.Range("I136").Value = 499.88635
$ws.Range("J136").Value = 1123.4546
$ws.Range("K136").Value = 1499.65905
$ws.Range("L136").Value = 3370.3638
$ws.Range("M136").Value = 1050.34095
$ws.Range("N136").Value = -8470.363799999999
